$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '37.376.73'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '2.067.24'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  +0.07%  '
Set-TextValue 'D5' '234.59'
$ws.Range('E5').Value = '  -0.92%  '
Set-TextValue 'D6' '0.625'
$ws.Range('E6').Value = '  +1.64%  '
$ws.Range('E7').Value = '  +0.04%  '
Set-TextValue 'D8' '57.38'
$ws.Range('E8').Value = '  -1.16%  '
$ws.Range('E9').Value = '  +3.51%  '
Set-TextValue 'D10' '0.0774'
$ws.Range('E10').Value = '  +1.49%  '
Set-TextValue 'D11' '0.103'
$ws.Range('E11').Value = '  +0.73%  '
$ws.Range('D12').Value = '2.373.14'
$ws.Range('E12').Value = '  +0.24%  '
Set-TextValue 'D13' '14.42'
$ws.Range('E13').Value = '  -0.29%  '
Set-TextValue 'D14' '20.70'
$ws.Range('E14').Value = '  -1.45%  '
Set-TextValue 'D15' '0.777'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('D17').Value = '2.071.33'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '37.335.32'
$ws.Range('E18').Value = '  -0.70%  '
Set-TextValue 'D19' '6.27'
$ws.Range('E19').Value = '  +1.51%  '
Set-TextValue 'D20' '69.64'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('D21').Value = '0.0₃0817'
$ws.Range('E21').Value = '  +0.16%  '
Set-TextValue 'D22' '226.76'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('E23').Value = '  -0.07%  '
Set-TextValue 'D25' '2.40'
$ws.Range('E25').Value = '  -2.00%  '
Set-TextValue 'D26' '166.96'
$ws.Range('E26').Value = '  +1.85%  '
Set-TextValue 'D27' '8.85'
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('E28').Value = '  -3.86%  '
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  -0.76%  '
Set-TextValue 'D32' '4.53'
$ws.Range('E32').Value = '  +0.14%  '
Set-TextValue 'D33' '0.0617'
$ws.Range('E33').Value = '  -0.78%  '
Set-TextValue 'D34' '4.54'
$ws.Range('E34').Value = '  +1.24%  '
Set-TextValue 'D35' '2.50'
$ws.Range('E35').Value = '  -3.26%  '
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('E37').Value = '  -3.02%  '
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('E39').Value = '  -4.54%  '
$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D40' '2.94'
$ws.Range('E40').Value = '  -0.72%  '
$ws.Range('B41').Value = 'Cronos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D41' '0.0964'
$ws.Range('E41').Value = '  -2.66%  '
Set-TextValue 'D42' '97.94'
$ws.Range('E42').Value = '  +0.59%  '
$ws.Range('D43').Value = '1.480.77'
$ws.Range('E43').Value = '  +0.27%  '
Set-TextValue 'D44' '0.0212'
$ws.Range('E44').Value = '  +0.75%  '
Set-TextValue 'D45' '1.16'
$ws.Range('E45').Value = '  -0.51%  '
Set-TextValue 'D46' '4.03'
$ws.Range('E46').Value = '  -12.31%  '
$ws.Range('E47').Value = '  -0.09%  '
Set-TextValue 'D48' '15.29'
$ws.Range('E48').Value = '  -3.78%  '
Set-TextValue 'D49' '7.24'
$ws.Range('E49').Value = '  +0.41%  '
Set-TextValue 'D50' '2.95'
$ws.Range('E50').Value = '  +0.61%  '
$ws.Range('D51').Value = '2.260.10'
$ws.Range('E51').Value = '  +0.22%  '
